$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Modbus")

$values = @{
    22 = "W"
    23 = "R"
    24 = "W"
    25 = "R"
    26 = "R"
    27 = "R"
    28 = "R"
    29 = "R"
    30 = "R"
    32 = "R"
    35 = "R"
    37 = "R"
    40 = "R"
    42 = "R"
    43 = "R"
    44 = "R"
    45 = "W"
    46 = "R"
    47 = "W"
}

$formatSource = $ws.Range("D2")
$formatSource.Copy()

foreach ($row in $values.Keys) {
    $cell = $ws.Range("D$row")
    $cell.PasteSpecial(-4122)
    $cell.Value = $values[$row]
}
